# Fruta / hortaliza, semanal
# A new weekly price-report row was inserted for "Feria Lagunitas de Puerto
# Montt" (Mandarina / Murcott / Primera) dated 44644, which pushes every
# subsequent record (previously rows 162-186) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 162 - this shifts rows 162:186
# down to 163:187 and grows the sheet dimension to A1:T187 automatically.
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row 162 with the new weekly record.
$ws.Range("A162").Value = 4
$ws.Range("B162").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C162").Value = "Los Lagos"
$ws.Range("D162").Value = 44644
$ws.Range("E162").Value = 10
$ws.Range("F162").Value = "Fruta"
$ws.Range("G162").Value = 100102
$ws.Range("H162").Value = "Cítricos"
$ws.Range("I162").Value = 100102004
$ws.Range("J162").Value = "Mandarina"
$ws.Range("K162").Value = "Murcott"
$ws.Range("L162").Value = "Primera"
$ws.Range("M162").Value = 400
$ws.Range("N162").Value = 12500
$ws.Range("O162").Value = 13000
$ws.Range("P162").Value = 12750
$ws.Range("Q162").Value = "$/bandeja 10 kilos"
$ws.Range("R162").Value = "Región de O'Higgins"
$ws.Range("S162").Value = 1275
$ws.Range("T162").Value = 10
